# Auto-generated PowerShell Excel COM-interop script
# Updates the '想去人数' (want-to-go count) values in column F
# across all four worksheets to match the gh-pages data refresh
# commit (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 42265
$ws.Range("F4").Value = 9972
$ws.Range("F5").Value = 222
$ws.Range("F6").Value = 1029
$ws.Range("F7").Value = 965
$ws.Range("F8").Value = 773
$ws.Range("F9").Value = 240
$ws.Range("F10").Value = 317
$ws.Range("F11").Value = 1008
$ws.Range("F14").Value = 798
$ws.Range("F15").Value = 346
$ws.Range("F16").Value = 1603
$ws.Range("F18").Value = 784
$ws.Range("F19").Value = 755
$ws.Range("F20").Value = 496
$ws.Range("F21").Value = 719
$ws.Range("F22").Value = 794
$ws.Range("F23").Value = 40
$ws.Range("F24").Value = 257
$ws.Range("F26").Value = 568
$ws.Range("F27").Value = 563
$ws.Range("F28").Value = 74
$ws.Range("F29").Value = 275
$ws.Range("F30").Value = 966
$ws.Range("F32").Value = 454
$ws.Range("F33").Value = 119
$ws.Range("F34").Value = 230
$ws.Range("F35").Value = 174
$ws.Range("F36").Value = 477
$ws.Range("F37").Value = 1410
$ws.Range("F38").Value = 322
$ws.Range("F39").Value = 1310
$ws.Range("F40").Value = 388
$ws.Range("F46").Value = 13
$ws.Range("F47").Value = 17

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 349
$ws.Range("F4").Value = 4479
$ws.Range("F8").Value = 153
$ws.Range("F10").Value = 64
$ws.Range("F16").Value = 35
$ws.Range("F17").Value = 4391

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 2100
$ws.Range("F4").Value = 469

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 2100
$ws.Range("F4").Value = 349
$ws.Range("F7").Value = 9972
$ws.Range("F8").Value = 222
$ws.Range("F9").Value = 1029
$ws.Range("F10").Value = 1029
$ws.Range("F12").Value = 965
$ws.Range("F13").Value = 153
$ws.Range("F14").Value = 240
$ws.Range("F15").Value = 317
$ws.Range("F16").Value = 1008
$ws.Range("F19").Value = 798
$ws.Range("F20").Value = 346
$ws.Range("F21").Value = 1603
$ws.Range("F23").Value = 755
$ws.Range("F24").Value = 496
$ws.Range("F25").Value = 719
$ws.Range("F26").Value = 794
$ws.Range("F27").Value = 40
$ws.Range("F29").Value = 568
$ws.Range("F32").Value = 563
$ws.Range("F33").Value = 74
$ws.Range("F34").Value = 275
$ws.Range("F37").Value = 454
$ws.Range("F38").Value = 119
$ws.Range("F39").Value = 230
$ws.Range("F41").Value = 1310
$ws.Range("F42").Value = 388
$ws.Range("F48").Value = 17
